# The "Good for how many people?" column (E) is blank for several menu
# items. Fill those in with "Not specified" so the data is explicit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (data starts at row 3; row 1 = title, row 2 = header) whose
# "Good for how many people?" cell is currently empty.
$blankRows = @(3, 5, 7, 8, 9, 10, 11, 13, 15, 17, 20, 21, 22, 23, 24, 25)

foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 5).Value = "Not specified"
}
